$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix C31: telephone was stored as text "33685678541.0" -> numeric 33685678541 ---
$ws.Range("C31").Value = 33685678541

# --- Insert a new row at 39 for "Misato Hayashi", pushing existing rows 39-49 down to 40-50 ---
$ws.Rows.Item(39).Insert()

$ws.Range("A39").Value = "Misato Hayashi"
$ws.Range("B39").Value = "Booking"

# telephone needs to stay textual (keeps the leading "+"), not get parsed as a number
$ws.Range("C39").NumberFormat = "@"
$ws.Range("C39").Value = "+14038525008"
$ws.Range("C39").Style = "Normal"

$ws.Range("D39").Value = 45900
$ws.Range("E39").Value = 45902
$ws.Range("F39").Value = 2
$ws.Range("G39").Value = 261.62
$ws.Range("H39").Value = 215.04
$ws.Range("I39").Value = 46.58
$ws.Range("J39").Value = 17.8
$ws.Range("K39").Value = 2025
$ws.Range("L39").Value = 8

# M/N/O stay blank for this row - clear any style/format inherited from the row insert
$ws.Range("M39").Style = "Normal"
$ws.Range("N39").Style = "Normal"
$ws.Range("O39").Style = "Normal"
